$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels (keep existing header style) ---
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# --- Data rows 2-11 ---
# Row 2: eng/POA
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "POA"
$ws.Range("C2").Value = "Proof of Address"
$ws.Range("D2").Value = "Address Proof"
$ws.Range("E2").Value = $true

# Row 3: fra/POA
$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "POA"
$ws.Range("C3").Value = "Un justificatif de domicile"
$ws.Range("D3").Value = "Preuve dadresse"
$ws.Range("E3").Value = $true

# Row 4: eng/POI
$ws.Range("A4").Value = "eng"
$ws.Range("B4").Value = "POI"
$ws.Range("C4").Value = "Proof of Identity"
$ws.Range("D4").Value = "Identity Proof"
$ws.Range("E4").Value = $true

# Row 5: fra/POI
$ws.Range("A5").Value = "fra"
$ws.Range("B5").Value = "POI"
$ws.Range("C5").Value = "Preuve didentité"
$ws.Range("D5").Value = "Preuve didentité"
$ws.Range("E5").Value = $true

# Row 6: eng/POR
$ws.Range("A6").Value = "eng"
$ws.Range("B6").Value = "POR"
$ws.Range("C6").Value = "Proof of Relationship"
$ws.Range("D6").Value = "Proof Relationship of the person"
$ws.Range("E6").Value = $true

# Row 7: fra/POR
$ws.Range("A7").Value = "fra"
$ws.Range("B7").Value = "POR"
$ws.Range("C7").Value = "Preuve de relation"
$ws.Range("D7").Value = "Preuve de relation de la personne"
$ws.Range("E7").Value = $true

# Row 8: eng/POB
$ws.Range("A8").Value = "eng"
$ws.Range("B8").Value = "POB"
$ws.Range("C8").Value = "Proof of Birth"
$ws.Range("D8").Value = "Proof date of birth of the person"
$ws.Range("E8").Value = $false

# Row 9: fra/POB
$ws.Range("A9").Value = "fra"
$ws.Range("B9").Value = "POB"
$ws.Range("C9").Value = "Preuve de naissance"
$ws.Range("D9").Value = "Preuve de la date de naissance de la personne"
$ws.Range("E9").Value = $false

# Row 10: eng/POE
$ws.Range("A10").Value = "eng"
$ws.Range("B10").Value = "POE"
$ws.Range("C10").Value = "Proof of Biometric Exception"
$ws.Range("D10").Value = "Proof of Biometric Exception"
$ws.Range("E10").Value = $true

# Row 11: fra/POE
$ws.Range("A11").Value = "fra"
$ws.Range("B11").Value = "POE"
$ws.Range("C11").Value = "Preuve dexception biométrique"
$ws.Range("D11").Value = "Preuve dexception biométrique"
$ws.Range("E11").Value = $true

# --- Column A (rows 2-11) uses the same bold/bordered style as the header row ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "done"
